$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.318.93"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.613.36"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.487"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.23%  "

$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").Value = "1.837.65"
$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").Value = "1.599.68"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "26.314.49"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.33%  "

$ws.Range("E21").Value = "  +1.23%  "

$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("E24").Value = "  +8.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.13%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -2.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("E30").Value = "  +3.80%  "

$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.41%  "

$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "1.161.14"
$ws.Range("E36").Value = "  +4.81%  "

$ws.Range("E37").Value = "  +9.17%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.502"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("E42").Value = "  +1.78%  "

$ws.Range("E43").Value = "  +3.00%  "

$ws.Range("D44").Value = "1.748.81"
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "

$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("D49").Value = "0.0₇0972"
$ws.Range("E49").Value = "  -14.41%  "

$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("E51").Value = "  -0.05%  "
